# Refactor for hidrometro - update quantity (F) and total value (H) figures
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H18").Value = 98903.18
$ws.Range("H22").Value = 6748.24
$ws.Range("H24").Value = 6970.97

$ws.Range("F36").Value = 1779.000
$ws.Range("H36").Value = 46339.13

$ws.Range("F38").Value = 320.000
$ws.Range("H38").Value = 27191.03

$ws.Range("F39").Value = 706.000
$ws.Range("H39").Value = 97247.97

$ws.Range("F45").Value = 2526.000
$ws.Range("H45").Value = 9764.55

$ws.Range("H70").Value = 35397.66
$ws.Range("H72").Value = 22265.59

$ws.Range("F91").Value = 9560.000
$ws.Range("H91").Value = 21698.08

$ws.Range("H92").Value = 938.55

$ws.Range("F99").Value = 1825.000
$ws.Range("H99").Value = 4581.39
